# Applies the crypto list refresh described in the target diff.
# D (Price) and E (Volume 1h %) columns are stored as literal text in the
# workbook (e.g. "59.270.69", "  +0.71%  "), so values are written with a
# leading quote-prefix ($q) to stop Excel from re-interpreting them as numbers
# and silently changing their formatting/precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$q = "'"  # forces text (quote-prefix) entry, mirrors the cells' original inline-string type

$ws.Range("D2").Value = $q + '59.270.69'

$ws.Range("D3").Value = $q + '2.586.31'
$ws.Range("E3").Value = $q + '  -0.42%  '

$ws.Range("E4").Value = $q + '  -0.01%  '

$ws.Range("D5").Value = $q + '571.38'
$ws.Range("E5").Value = $q + '  +3.57%  '

$ws.Range("D6").Value = $q + '143.59'
$ws.Range("E6").Value = $q + '  +0.29%  '

$ws.Range("D7").Value = $q + '0.999'
$ws.Range("E7").Value = $q + '  +0.04%  '

$ws.Range("D8").Value = $q + '0.602'
$ws.Range("E8").Value = $q + '  -0.66%  '

$ws.Range("D9").Value = $q + '2.596.67'
$ws.Range("E9").Value = $q + '  -0.45%  '

$ws.Range("E10").Value = $q + '  -1.70%  '

$ws.Range("E11").Value = $q + '  +3.06%  '

$ws.Range("D12").Value = $q + '0.156'
$ws.Range("E12").Value = $q + '  +9.80%  '

$ws.Range("E13").Value = $q + '  +2.87%  '

$ws.Range("D14").Value = $q + '3.042.14'
$ws.Range("E14").Value = $q + '  -0.41%  '

$ws.Range("D15").Value = $q + '59.311.02'

$ws.Range("D16").Value = $q + '22.58'
$ws.Range("E16").Value = $q + '  +8.10%  '

$ws.Range("E17").Value = $q + '  +4.10%  '

$ws.Range("D18").Value = $q + '2.590.84'
$ws.Range("E18").Value = $q + '  -0.25%  '

$ws.Range("D19").Value = $q + '4.53'
$ws.Range("E19").Value = $q + '  +1.30%  '

$ws.Range("D20").Value = $q + '337.84'
$ws.Range("E20").Value = $q + '  +0.15%  '

$ws.Range("D21").Value = $q + '10.23'
$ws.Range("E21").Value = $q + '  +1.71%  '

$ws.Range("E22").Value = $q + '  +0.52%  '

$ws.Range("E23").Value = $q + '  +0.05%  '

$ws.Range("D24").Value = $q + '64.35'
$ws.Range("E24").Value = $q + '  -3.93%  '

$ws.Range("E25").Value = $q + '  +5.86%  '

$ws.Range("D26").Value = $q + '1.00'
$ws.Range("E26").Value = $q + '  +0.19%  '

$ws.Range("E27").Value = $q + '  +1.65%  '

$ws.Range("D28").Value = $q + '7.29'
$ws.Range("E28").Value = $q + '  +1.94%  '

$ws.Range("D29").Value = $q + '0.0₃0781'
$ws.Range("E29").Value = $q + '  +3.75%  '

$ws.Range("E30").Value = $q + '  +0.01%  '

$ws.Range("E31").Value = $q + '  +0.85%  '

$ws.Range("D32").Value = $q + '6.08'
$ws.Range("E32").Value = $q + '  +1.47%  '

$ws.Range("D33").Value = $q + '158.31'
$ws.Range("E33").Value = $q + '  +2.31%  '

$ws.Range("D34").Value = $q + '19.07'
$ws.Range("E34").Value = $q + '  +0.55%  '

$ws.Range("E35").Value = $q + '  +2.60%  '

$ws.Range("E36").Value = $q + '  +1.96%  '

$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").Value = $q + '0.885'
$ws.Range("E37").Value = $q + '  +6.81%  '

$ws.Range("B38").Value = 'SuiNetwork'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D38").Value = $q + '0.880'
$ws.Range("E38").Value = $q + '  -1.30%  '

$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = $q + '37.04'
$ws.Range("E39").Value = $q + '  +0.13%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = $q + '1.50'
$ws.Range("E40").Value = $q + '  +2.60%  '

$ws.Range("D41").Value = $q + '295.01'
$ws.Range("E41").Value = $q + '  +3.94%  '

$ws.Range("E42").Value = $q + '  +1.73%  '

$ws.Range("E43").Value = $q + '  +0.03%  '

$ws.Range("E44").Value = $q + '  +1.90%  '

$ws.Range("D45").Value = $q + '0.596'
$ws.Range("E45").Value = $q + '  -0.47%  '

$ws.Range("D46").Value = $q + '0.0537'
$ws.Range("E46").Value = $q + '  +0.69%  '

$ws.Range("E47").Value = $q + '  +2.90%  '

$ws.Range("E48").Value = $q + '  +0.03%  '

$ws.Range("D49").Value = $q + '124.48'
$ws.Range("E49").Value = $q + '  +5.03%  '

$ws.Range("E50").Value = $q + '  +2.31%  '

$ws.Range("D51").Value = $q + '1.946.14'
$ws.Range("E51").Value = $q + '  -0.02%  '
